$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values for rows 2-5 (column A "cluster id", column B "count")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 392

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 248

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 218

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 143

# Remove the old row 6 entirely (data + shrinks used range/dimension)
$ws.Range("A6:B6").Delete()
